# Auto-generated edit script: refresh cached market-data values
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) per the
# scheduled-runner data refresh described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 27785278
$ws.Range("I62").Value = 27785278
$ws.Range("K62").Value = 27785278
$ws.Range("M62").Value = -27784654
$ws.Range("H65").Value = 27785278
$ws.Range("I65").Value = 27785278
$ws.Range("K65").Value = 138926390
$ws.Range("M65").Value = -138923270
$ws.Range("H70").Value = 1715.5555
$ws.Range("I70").Value = 1640
$ws.Range("J70").Value = 1776
$ws.Range("K70").Value = 4920
$ws.Range("L70").Value = 5328
$ws.Range("M70").Value = -4650
$ws.Range("N70").Value = -5868
$ws.Range("H73").Value = 1715.5555
$ws.Range("I73").Value = 1640
$ws.Range("J73").Value = 1776
$ws.Range("K73").Value = 4920
$ws.Range("L73").Value = 5328
$ws.Range("M73").Value = -3984
$ws.Range("N73").Value = -7200
$ws.Range("H125").Value = 2265.3333
$ws.Range("I125").Value = 394
$ws.Range("J125").Value = 2800
$ws.Range("K125").Value = 3546
$ws.Range("L125").Value = 25200
$ws.Range("M125").Value = -1086
$ws.Range("N125").Value = -30120
$ws.Range("H135").Value = 32258644
$ws.Range("I135").Value = 355.39285
$ws.Range("J135").Value = 333336000
$ws.Range("K135").Value = 3198.53565
$ws.Range("L135").Value = 3000024000
$ws.Range("M135").Value = -663.5356500000003
$ws.Range("N135").Value = -3000029070
$ws.Range("H137").Value = 1598.4193
$ws.Range("I137").Value = 1120
$ws.Range("J137").Value = 1900.579
$ws.Range("K137").Value = 3360
$ws.Range("L137").Value = 5701.737
$ws.Range("M137").Value = -810
$ws.Range("N137").Value = -10801.737

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 3600
$ws.Range("J38").Value = 3600
$ws.Range("L38").Value = 3600
$ws.Range("N38").Value = -4432
$ws.Range("H39").Value = 1500
$ws.Range("J39").Value = 1500
$ws.Range("L39").Value = 1500
$ws.Range("N39").Value = -2278
$ws.Range("H44").Value = 8000
$ws.Range("J44").Value = 8000
$ws.Range("L44").Value = 8000
$ws.Range("N44").Value = -8994
$ws.Range("H46").Value = 6500
$ws.Range("J46").Value = 6500
$ws.Range("L46").Value = 6500
$ws.Range("N46").Value = -7096
$ws.Range("H49").Value = 10000
$ws.Range("J49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("N49").Value = -10478
$ws.Range("H132").Value = 1458102.1
$ws.Range("J132").Value = 1458102.1
$ws.Range("L132").Value = 1458102.1
$ws.Range("N132").Value = -1468222.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 689.1667
$ws.Range("I122").Value = 695.5
$ws.Range("K122").Value = 2086.5
$ws.Range("M122").Value = 363.5
$ws.Range("H132").Value = 1824.5217
$ws.Range("I132").Value = 1298.4286
$ws.Range("J132").Value = 2642.889
$ws.Range("K132").Value = 3895.2858
$ws.Range("L132").Value = 7928.667
$ws.Range("M132").Value = -1365.2858
$ws.Range("N132").Value = -12988.667
$ws.Range("H134").Value = 15626628
$ws.Range("I134").Value = 1688.6923
$ws.Range("J134").Value = 83334696
$ws.Range("K134").Value = 5066.0769
$ws.Range("L134").Value = 250004088
$ws.Range("M134").Value = -2531.0769
$ws.Range("N134").Value = -250009158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 433.33334
$ws.Range("I15").Value = 433.33334
$ws.Range("K15").Value = 1300.00002
$ws.Range("M15").Value = -1160.00002
$ws.Range("H60").Value = 2000.8334
$ws.Range("I60").Value = 752.5
$ws.Range("J60").Value = 2250.5
$ws.Range("K60").Value = 2257.5
$ws.Range("L60").Value = 6751.5
$ws.Range("M60").Value = -2006.5
$ws.Range("N60").Value = -7253.5
$ws.Range("H81").Value = 3466.3333
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 3699.5
$ws.Range("K81").Value = 9000
$ws.Range("L81").Value = 11098.5
$ws.Range("M81").Value = -7877
$ws.Range("N81").Value = -13344.5
$ws.Range("H84").Value = 3466.3333
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 3699.5
$ws.Range("K84").Value = 27000
$ws.Range("L84").Value = 33295.5
$ws.Range("M84").Value = -21384
$ws.Range("N84").Value = -44527.5
$ws.Range("H132").Value = 698.7692
$ws.Range("I132").Value = 698.7692
$ws.Range("K132").Value = 6288.922799999999
$ws.Range("M132").Value = -3758.922799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 263.54544
$ws.Range("I2").Value = 250.33333
$ws.Range("J2").Value = 279.4
$ws.Range("K2").Value = 250.33333
$ws.Range("L2").Value = 279.4
$ws.Range("M2").Value = -137.33333
$ws.Range("N2").Value = -505.4
$ws.Range("H21").Value = 1669800
$ws.Range("I21").Value = 5000000
$ws.Range("J21").Value = 4700
$ws.Range("K21").Value = 5000000
$ws.Range("L21").Value = 4700
$ws.Range("M21").Value = -4999827
$ws.Range("N21").Value = -5046
$ws.Range("H30").Value = 1669800
$ws.Range("I30").Value = 5000000
$ws.Range("J30").Value = 4700
$ws.Range("K30").Value = 5000000
$ws.Range("L30").Value = 4700
$ws.Range("M30").Value = -4999895
$ws.Range("N30").Value = -4910
$ws.Range("H122").Value = 8066326
$ws.Range("I122").Value = 1914.1
$ws.Range("J122").Value = 22728894
$ws.Range("K122").Value = 5742.299999999999
$ws.Range("L122").Value = 68186682
$ws.Range("M122").Value = -3292.299999999999
$ws.Range("N122").Value = -68191582
$ws.Range("H132").Value = 2303.0286
$ws.Range("I132").Value = 2255.652
$ws.Range("K132").Value = 6766.956
$ws.Range("M132").Value = -4236.956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 715.61536
$ws.Range("I22").Value = 476.16666
$ws.Range("K22").Value = 476.16666
$ws.Range("M22").Value = -181.16666
$ws.Range("H27").Value = 715.61536
$ws.Range("I27").Value = 476.16666
$ws.Range("K27").Value = 476.16666
$ws.Range("M27").Value = -369.16666
$ws.Range("H100").Value = 1343
$ws.Range("J100").Value = 1371
$ws.Range("L100").Value = 1371
$ws.Range("N100").Value = -2453
$ws.Range("H136").Value = 1511.8889
$ws.Range("I136").Value = 1325.875
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3977.625
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1427.625
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178
$ws.Range("H132").Value = 2600.9512
$ws.Range("I132").Value = 2716
$ws.Range("K132").Value = 8148
$ws.Range("M132").Value = -5618
